$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "started at midnight" note that was sitting in D38.
$ws.Range("D38").ClearContents()

# The AI-script recording session actually ran 1.5 hours, not 2.
$ws.Range("B37").Value = 1.5

# New row: got the simple AI script working in the game.
$ws.Range("C39").Copy() | Out-Null
$ws.Range("C40").PasteSpecial(-4122) | Out-Null
$ws.Range("A40").Value = "Getting Simple Ai Script into the game."
$ws.Range("B40").Value = 4
$ws.Range("C40").Value = 41926

$ws.Range("I40").Select() | Out-Null
